# Update the build timestamp embedded in the version string, wherever it
# appears in the workbook ("About" sheet text and the per-row version
# column on the "Boundaries and methane sources" sheet).

$wb = $excel.ActiveWorkbook

$oldStamp = "February 03 2026 17.29.55 EST"
$newStamp = "February 03 2026 18.05.36 EST"

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        $val = $cell.Value()
        if ($val -ne $null -and $val -is [string] -and $val.Contains($oldStamp)) {
            $cell.Value = $val.Replace($oldStamp, $newStamp)
        }
    }
}
